$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the old "Hold type" / "Hold color" / "[List of] routes" headers
# (D1:F1) one column to the right (E1:G1), right-to-left so nothing is
# clobbered before it's read. Done as plain value writes (not a structural
# column insert) so column widths/formatting stay exactly as they were.
$ws.Range("G1").Value = $ws.Range("F1").Text
$ws.Range("F1").Value = $ws.Range("E1").Text
$ws.Range("E1").Value = $ws.Range("D1").Text

# Replace "Hold pos. X" / "Hold pos. Y" with the new labels, and add the
# new "spot-pos" column in D1.
$ws.Range("B1").Value = "X-loc"
$ws.Range("C1").Value = "Y-loc"
$ws.Range("D1").Value = "spot-pos"

# Match the saved selection state
$ws.Range("D1").Select()
